# Apply the renames + view/selection changes described by the commit:
#   "I've done loading and preproc function. Trouble with sheets names is
#    also allegedly rid off using RegEx."
#
# 1) Shorten the four sheet names (RegEx clean-up of the verbose
#    "input_..."/"..._names" naming).
# 2) Update each sheet's remembered selection / active cell.
# 3) Leave the last-touched sheet (component) as the active tab, which
#    also drives workbookView's activeTab and the per-sheet
#    tabSelected flag.

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets (stripped via regex-like shortening) ---------------
$wb.Worksheets.Item(1).Name = "stoich_coeff"
$wb.Worksheets.Item(2).Name = "constants_log10"
$wb.Worksheets.Item(3).Name = "concentra"
$wb.Worksheets.Item(4).Name = "component"

# --- 2. Update selections on each affected sheet --------------------------
# stoich_coeff: J11 -> G31
$wsStoich = $wb.Worksheets.Item(1)
$wsStoich.Activate() | Out-Null
$wsStoich.Range("G31").Select() | Out-Null

# concentra: F2:F13 -> H8 (also loses tabSelected as focus moves on)
$wsConcentra = $wb.Worksheets.Item(3)
$wsConcentra.Activate() | Out-Null
$wsConcentra.Range("H8").Select() | Out-Null

# component: gains tabSelected + new selection I18, and ends as active tab
$wsComponent = $wb.Worksheets.Item(4)
$wsComponent.Activate() | Out-Null
$wsComponent.Range("I18").Select() | Out-Null
